$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 gains a value in column A (new shared string) and its row height grows
# from 21.6 to 43.2, matching the other single-column "script filename" rows.
$ws.Range("A21").Value = "SCRIPT/P02P01A/us0201.ssb"

# New rows 22-25, each holding a single script-filename value in column A,
# following the same pattern/style as rows 4-14 and 22-25 in the final sheet.
$ws.Range("A22").Value = "SCRIPT/P02P01A/us0401.ssb"
$ws.Range("A23").Value = "SCRIPT/P02P01A/us2001.ssb"
$ws.Range("A24").Value = "SCRIPT/P02P01A/us2004.ssb"
$ws.Range("A25").Value = "SCRIPT/P02P01A/us2007.ssb"

# Match the 43.2pt row height used by every other "filename-only" row.
$ws.Range("A21:A25").RowHeight = 43.2

# Move the selection to the new last cell, like the source workbook after editing.
[void]$ws.Range("E21").Select()
